# Generate Report for Handback
# Adds a new handback record (b4821fdf-1350-4548-b609-44b349cd1daf.md) as
# row 4 on the "Overview", "zh-cn" and "de-de" sheets, keeping each
# sheet's Excel Table (ListObject) in sync and re-creating the matching
# hyperlinks.

$wb = $excel.ActiveWorkbook

$fileName = "b4821fdf-1350-4548-b609-44b349cd1daf.md"
$pathAndName = "e2e\b4821fdf-1350-4548-b609-44b349cd1daf.md"
$status = "Handed back: in sync with en-US"
$handoffDate = "2016-08-12 20:54:34"

$zhXlf = "b4821fdf-1350-4548-b609-44b349cd1daf.6d027c4b7d9cad09bc8c07834cf4555aadc20cc6.zh-cn.xlf"
$zhHandoffDate = "2016-08-12 20:54:27"
$zhHandbackDate = "2016-08-12 20:54:56"

$deXlf = "b4821fdf-1350-4548-b609-44b349cd1daf.6d027c4b7d9cad09bc8c07834cf4555aadc20cc6.de-de.xlf"
$deHandoffDate = "2016-08-12 20:55:12"

$ghBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/99ad2b511c3bc47826f6de9f6b15587b15ffa35a/e2e/$fileName"
$ghZh = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/752aafcac0ecee83f4d15a7e9a0a3450cb075858/e2e/$fileName"
$ghDe = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/521ed2a21329da14ca9b95a1571e378ae2dbc6aa/e2e/$fileName"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item("Overview")
$tblOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("B4").Value = $pathAndName
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $status
$wsOverview.Range("F4").Value = $status
$wsOverview.Range("G4").Value = $handoffDate
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $ghBase, "", "", $pathAndName) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$tblZh = $wsZh.ListObjects.Item("zh-cn")
$tblZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = $fileName
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = $zhHandoffDate
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = $fileName
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = $zhHandbackDate
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $ghBase, "", "", $fileName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $ghZh, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$tblDe = $wsDe.ListObjects.Item("de-de")
$tblDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = $fileName
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = $handoffDate
$wsDe.Range("I4").Value = $fileName
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = $deHandoffDate
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $ghBase, "", "", $fileName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $ghDe, "", "", $fileName) | Out-Null
